$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1083.0266
$ws.Range("I15").Value = 1083.0266
$ws.Range("K15").Value = 3249.0798
$ws.Range("M15").Value = -3080.0798

$ws.Range("H40").Value = 2215.6428
$ws.Range("I40").Value = 2191.6667
$ws.Range("K40").Value = 2191.6667
$ws.Range("M40").Value = -2016.6667

$ws.Range("H51").Value = 17171.715
$ws.Range("J51").Value = 4000.6667
$ws.Range("L51").Value = 4000.6667
$ws.Range("N51").Value = -4968.6667

$ws.Range("H61").Value = 5589938.5
$ws.Range("I61").Value = 6944919
$ws.Range("K61").Value = 20834757
$ws.Range("M61").Value = -20834585

$ws.Range("H98").Value = 1095.0416
$ws.Range("I98").Value = 1017.8823
$ws.Range("J98").Value = 1282.4286
$ws.Range("K98").Value = 1017.8823
$ws.Range("L98").Value = 1282.4286
$ws.Range("M98").Value = 480.1177
$ws.Range("N98").Value = -4278.4286

$ws.Range("H106").Value = 2541.2856
$ws.Range("I106").Value = 2541.2856
$ws.Range("K106").Value = 2541.2856
$ws.Range("M106").Value = -1910.2856

$ws.Range("H122").Value = 1095.0416
$ws.Range("I122").Value = 1017.8823
$ws.Range("J122").Value = 1282.4286
$ws.Range("K122").Value = 3053.6469
$ws.Range("L122").Value = 3847.2858
$ws.Range("M122").Value = -603.6468999999997
$ws.Range("N122").Value = -8747.2858

$ws.Range("H135").Value = 665.0540999999999
$ws.Range("I135").Value = 643.48486
$ws.Range("J135").Value = 843
$ws.Range("K135").Value = 5791.363740000001
$ws.Range("L135").Value = 7587
$ws.Range("M135").Value = -3256.363740000001
$ws.Range("N135").Value = -12657

$ws.Range("H137").Value = 2379.8667
$ws.Range("I137").Value = 1443.1666
$ws.Range("J137").Value = 3004.3333
$ws.Range("K137").Value = 4329.4998
$ws.Range("L137").Value = 9012.999899999999
$ws.Range("M137").Value = -1779.4998
$ws.Range("N137").Value = -14112.9999

$ws.Range("H138").Value = 6067.032
$ws.Range("I138").Value = 1786.4736
$ws.Range("J138").Value = 12844.583
$ws.Range("K138").Value = 5359.4208
$ws.Range("L138").Value = 38533.749
$ws.Range("M138").Value = -219.4207999999999
$ws.Range("N138").Value = -48813.749

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24269.4
$ws.Range("I32").Value = 4158.797
$ws.Range("K32").Value = 4158.797
$ws.Range("M32").Value = -3871.797

$ws.Range("H61").Value = 1802.5962
$ws.Range("I61").Value = 990.37036
$ws.Range("J61").Value = 2679.8
$ws.Range("K61").Value = 990.37036
$ws.Range("L61").Value = 2679.8
$ws.Range("M61").Value = -778.37036
$ws.Range("N61").Value = -3103.8

$ws.Range("H74").Value = 1641.909
$ws.Range("I74").Value = 1920.4615
$ws.Range("J74").Value = 1525.0968
$ws.Range("K74").Value = 1920.4615
$ws.Range("L74").Value = 1525.0968
$ws.Range("M74").Value = -1046.4615
$ws.Range("N74").Value = -3273.0968

$ws.Range("H77").Value = 1641.909
$ws.Range("I77").Value = 1920.4615
$ws.Range("J77").Value = 1525.0968
$ws.Range("K77").Value = 9602.307499999999
$ws.Range("L77").Value = 7625.484
$ws.Range("M77").Value = -5234.307499999999
$ws.Range("N77").Value = -16361.484

$ws.Range("H132").Value = 2543.1091
$ws.Range("I132").Value = 2627.4773
$ws.Range("J132").Value = 2205.6365
$ws.Range("K132").Value = 7882.4319
$ws.Range("L132").Value = 6616.9095
$ws.Range("M132").Value = -5352.4319
$ws.Range("N132").Value = -11676.9095

$ws.Range("H136").Value = 1802.5962
$ws.Range("I136").Value = 990.37036
$ws.Range("J136").Value = 2679.8
$ws.Range("K136").Value = 2971.11108
$ws.Range("L136").Value = 8039.400000000001
$ws.Range("M136").Value = -421.1110800000001
$ws.Range("N136").Value = -13139.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1966.7222
$ws.Range("I20").Value = 1580.1428
$ws.Range("J20").Value = 3319.75
$ws.Range("K20").Value = 1580.1428
$ws.Range("L20").Value = 3319.75
$ws.Range("M20").Value = -1333.1428
$ws.Range("N20").Value = -3813.75

$ws.Range("H57").Value = 53000
$ws.Range("J57").Value = 68000
$ws.Range("L57").Value = 68000
$ws.Range("N57").Value = -69440

$ws.Range("H86").Value = 15560.875
$ws.Range("J86").Value = 3871.75
$ws.Range("L86").Value = 3871.75
$ws.Range("N86").Value = -6117.75

$ws.Range("H89").Value = 15560.875
$ws.Range("J89").Value = 3871.75
$ws.Range("L89").Value = 19358.75
$ws.Range("N89").Value = -30590.75

$ws.Range("H99").Value = 1687.7142
$ws.Range("I99").Value = 1279.75
$ws.Range("K99").Value = 1279.75
$ws.Range("M99").Value = 218.25

$ws.Range("H134").Value = 1290.228
$ws.Range("I134").Value = 1332.289
$ws.Range("J134").Value = 1132.5
$ws.Range("K134").Value = 3996.867
$ws.Range("L134").Value = 3397.5
$ws.Range("M134").Value = -1461.867
$ws.Range("N134").Value = -8467.5

$ws.Range("H136").Value = 53000
$ws.Range("J136").Value = 68000
$ws.Range("L136").Value = 68000
$ws.Range("N136").Value = -78200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18342.7
$ws.Range("I31").Value = 44606.13
$ws.Range("J31").Value = 2016.7838
$ws.Range("K31").Value = 44606.13
$ws.Range("L31").Value = 2016.7838
$ws.Range("M31").Value = -44311.13
$ws.Range("N31").Value = -2606.7838

$ws.Range("H34").Value = 18342.7
$ws.Range("I34").Value = 44606.13
$ws.Range("J34").Value = 2016.7838
$ws.Range("K34").Value = 44606.13
$ws.Range("L34").Value = 2016.7838
$ws.Range("M34").Value = -44404.13
$ws.Range("N34").Value = -2420.7838

$ws.Range("H51").Value = 6992.25
$ws.Range("I51").Value = 4070
$ws.Range("J51").Value = 7966.3335
$ws.Range("K51").Value = 4070
$ws.Range("L51").Value = 7966.3335
$ws.Range("M51").Value = -3334
$ws.Range("N51").Value = -9438.333500000001

$ws.Range("H61").Value = 6992.25
$ws.Range("I61").Value = 4070
$ws.Range("J61").Value = 7966.3335
$ws.Range("K61").Value = 4070
$ws.Range("L61").Value = 7966.3335
$ws.Range("M61").Value = -3722
$ws.Range("N61").Value = -8662.333500000001

$ws.Range("H86").Value = 2531.56
$ws.Range("I86").Value = 1722.3077
$ws.Range("J86").Value = 3408.25
$ws.Range("K86").Value = 1722.3077
$ws.Range("L86").Value = 3408.25
$ws.Range("M86").Value = -599.3077000000001
$ws.Range("N86").Value = -5654.25

$ws.Range("H89").Value = 2531.56
$ws.Range("I89").Value = 1722.3077
$ws.Range("J89").Value = 3408.25
$ws.Range("K89").Value = 8611.538500000001
$ws.Range("L89").Value = 17041.25
$ws.Range("M89").Value = -2995.538500000001
$ws.Range("N89").Value = -28273.25

$ws.Range("H93").Value = 19042.666
$ws.Range("I93").Value = 9951.75
$ws.Range("J93").Value = 26315.4
$ws.Range("K93").Value = 9951.75
$ws.Range("L93").Value = 26315.4
$ws.Range("M93").Value = -8079.75
$ws.Range("N93").Value = -30059.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 18760.186
$ws.Range("I68").Value = 1122.6923
$ws.Range("J68").Value = 23744.695
$ws.Range("K68").Value = 3368.0769
$ws.Range("L68").Value = 71234.08499999999
$ws.Range("M68").Value = -2557.0769
$ws.Range("N68").Value = -72856.08499999999

$ws.Range("H71").Value = 18760.186
$ws.Range("I71").Value = 1122.6923
$ws.Range("J71").Value = 23744.695
$ws.Range("K71").Value = 10104.2307
$ws.Range("L71").Value = 213702.255
$ws.Range("M71").Value = -6048.2307
$ws.Range("N71").Value = -221814.255

$ws.Range("H103").Value = 124.4
$ws.Range("I103").Value = 124.4
$ws.Range("K103").Value = 373.2
$ws.Range("M103").Value = 505.8

$ws.Range("H131").Value = 895.8163500000001
$ws.Range("I131").Value = 676
$ws.Range("J131").Value = 902.7578999999999
$ws.Range("K131").Value = 2028
$ws.Range("L131").Value = 2708.2737
$ws.Range("M131").Value = 3012
$ws.Range("N131").Value = -12788.2737

$ws.Range("H133").Value = 2745
$ws.Range("I133").Value = 1137.1428
$ws.Range("J133").Value = 14000
$ws.Range("K133").Value = 3411.4284
$ws.Range("L133").Value = 42000
$ws.Range("M133").Value = 1648.5716
$ws.Range("N133").Value = -52120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H132").Value = 2354.0833
$ws.Range("I132").Value = 2056.7334
$ws.Range("K132").Value = 6170.2002
$ws.Range("M132").Value = -3640.2002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5195.6816
$ws.Range("I132").Value = 5614.857
$ws.Range("K132").Value = 16844.571
$ws.Range("M132").Value = -14314.571

$ws.Range("H136").Value = 1831.4375
$ws.Range("I136").Value = 1577.7778
$ws.Range("J136").Value = 2157.5715
$ws.Range("K136").Value = 4733.3334
$ws.Range("L136").Value = 6472.7145
$ws.Range("M136").Value = -2183.3334
$ws.Range("N136").Value = -11572.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 16500
$ws.Range("I58").Value = 16500
$ws.Range("K58").Value = 16500
$ws.Range("M58").Value = -16192

$ws.Range("H136").Value = 786.2941
$ws.Range("I136").Value = 556.0476
$ws.Range("K136").Value = 1668.1428
$ws.Range("M136").Value = 881.8571999999999
